# Changes of DEV URL configuration
# Update ShipmentTrackNum (column C) values on the first sheet with a new
# batch of tracking numbers; rows that mirrored the value into column D
# (PackageTrackNum) get the same new value there too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTrackNums = @{
    2  = "320017961973"
    3  = "320017961984"
    4  = "320017962010"
    5  = "320017962031"
    6  = "320017962075"
    7  = "320017962097"
    8  = "320017962123"
    9  = "320017962145"
    10 = "320017962178"
    11 = "320017962190"
    12 = "320017962237"
    13 = "320017962259"
    14 = "320017962281"
    15 = "320017962307"
    16 = "320017962330"
    17 = "320017962351"
    18 = "320017962395"
    19 = "320017962410"
    20 = "320017962443"
    21 = "320017962465"
    22 = "320017962498"
}

# Rows whose column D mirrors column C's tracking number.
$mirrorRows = @(5, 6, 7, 13, 14, 15, 16, 17)

foreach ($row in $newTrackNums.Keys) {
    $value = $newTrackNums[$row]
    $ws.Cells.Item($row, 3) = $value
    if ($mirrorRows -contains $row) {
        $ws.Cells.Item($row, 4) = $value
    }
}
